{"js": "// Apply the \"Add cellphone texts, and stairs to greenhouse\" edit.\n// Six textual tweaks to Diego's morning-routine dialogue, plus moving the\n// \"_GoBack\" bookmark from the end of the \"Diego\" line to wrap \"Let's go\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: find the (single) paragraph whose text matches `predicate`.\nfunction findParagraph(predicate) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (predicate(paragraphs.items[i].text)) return paragraphs.items[i];\n  }\n  throw new Error(\"paragraph not found\");\n}\n\n// --- 1) \"D- I need cook something, shower, and find something nice to wear.\"\n//        -> \"D- I need to cook something, shower, and find something nice to wear.\"\nconst pNeedNice = findParagraph(\n  (t) => t === \"D- I need cook something, shower, and find something nice to wear.\"\n);\n{\n  const range = pNeedNice.getRange();\n  const hits = range.search(\"D- I need \", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertText(\"D- I need to \", \"Replace\");\n  await context.sync();\n}\n\n// --- 2) \"D- I need cook something, shower, and find something to wear.\"\n//        -> \"D- I need to cook something, shower, and find something to wear.\"\nconst pNeedPlain = findParagraph(\n  (t) => t === \"D- I need cook something, shower, and find something to wear.\"\n);\n{\n  const range = pNeedPlain.getRange();\n  const hits = range.search(\"D- I need \", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertText(\"D- I need to \", \"Replace\");\n  await context.sync();\n}\n\n// --- 3) \"D-I don't want to think about it, it stresses me out\" -> add trailing \".\"\nconst pStress = findParagraph(\n  (t) => t === \"D-I don\\u2019t want to think about it, it stresses me out\"\n);\n{\n  const range = pStress.getRange();\n  range.insertText(\".\", \"End\");\n  await context.sync();\n}\n\n// --- 4) \"D- Great, I'm late. I'll just grab whatever clothes I have laying around\"\n//        -> \"...grab whichever clothes I have laying around.\"\nconst pGrab = findParagraph(\n  (t) =>\n    t ===\n    \"D- Great, I\\u2019m late. I\\u2019ll just grab whatever clothes I have laying around\"\n);\n{\n  const range = pGrab.getRange();\n  const hits = range.search(\"whatever\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertText(\"whichever\", \"Replace\");\n  await context.sync();\n  range.insertText(\".\", \"End\");\n  await context.sync();\n}\n\n// --- 5) \"*sniff sniff* Ok it's not too bad\" -> add trailing \".\"\nconst pSniff = findParagraph((t) => t === \"*sniff sniff* Ok it\\u2019s not too bad\");\n{\n  const range = pSniff.getRange();\n  range.insertText(\".\", \"End\");\n  await context.sync();\n}\n\n// --- 6) Move the \"_GoBack\" bookmark: remove it from after \"Diego\" and wrap\n//        \"Let's go\" (in the \"D- Let's go\" line) with it instead.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst pGo = findParagraph((t) => t === \"D- Let\\u2019s go\");\n{\n  const range = pGo.getRange();\n  const hits = range.search(\"Let\\u2019s go\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Apply the \"Add cellphone texts, and stairs to greenhouse\" edit.\n# Six textual tweaks to Diego's morning-routine dialogue, plus moving the\n# \"_GoBack\" bookmark from the end of the \"Diego\" line to wrap \"Let's go\".\n\n$d = $word.ActiveDocument\n$RS = [char]0x2019   # right single quotation mark used throughout the script\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    # --- 1) \"D- I need cook something, shower, and find something nice to wear.\"\n    #        -> \"D- I need to cook something, shower, and find something nice to wear.\"\n    if ($t -eq (\"D- I need cook something, shower, and find something nice to wear.`r\")) {\n        $r = $p.Range.Duplicate\n        $find = $r.Find\n        $find.Text = \"D- I need \"\n        $find.Execute() | Out-Null\n        $r.Text = \"D- I need to \"\n    }\n\n    # --- 2) \"D- I need cook something, shower, and find something to wear.\"\n    #        -> \"D- I need to cook something, shower, and find something to wear.\"\n    elseif ($t -eq (\"D- I need cook something, shower, and find something to wear.`r\")) {\n        $r = $p.Range.Duplicate\n        $find = $r.Find\n        $find.Text = \"D- I need \"\n        $find.Execute() | Out-Null\n        $r.Text = \"D- I need to \"\n    }\n\n    # --- 3) \"D-I don't want to think about it, it stresses me out\" -> add trailing \".\"\n    elseif ($t -eq (\"D-I don\" + $RS + \"t want to think about it, it stresses me out`r\")) {\n        $r = $p.Range.Duplicate\n        $r.MoveEnd(1, -1)\n        $r.InsertAfter(\".\")\n    }\n\n    # --- 4) \"D- Great, I'm late. I'll just grab whatever clothes I have laying around\"\n    #        -> \"...grab whichever clothes I have laying around.\"\n    elseif ($t -eq (\"D- Great, I\" + $RS + \"m late. I\" + $RS + \"ll just grab whatever clothes I have laying around`r\")) {\n        $r = $p.Range.Duplicate\n        $find = $r.Find\n        $find.Text = \"whatever\"\n        $find.Execute() | Out-Null\n        $r.Text = \"whichever\"\n\n        $r2 = $p.Range.Duplicate\n        $r2.MoveEnd(1, -1)\n        $r2.InsertAfter(\".\")\n    }\n\n    # --- 5) \"*sniff sniff* Ok it's not too bad\" -> add trailing \".\"\n    elseif ($t -eq (\"*sniff sniff* Ok it\" + $RS + \"s not too bad`r\")) {\n        $r = $p.Range.Duplicate\n        $r.MoveEnd(1, -1)\n        $r.InsertAfter(\".\")\n    }\n}\n\n# --- 6) Move the \"_GoBack\" bookmark: remove it from after \"Diego\" and wrap\n#        \"Let's go\" (in the \"D- Let's go\" line) with it instead.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -eq (\"D- Let\" + $RS + \"s go`r\")) {\n        $r = $p.Range.Duplicate\n        $find = $r.Find\n        $find.Text = \"Let\" + $RS + \"s go\"\n        $find.Execute() | Out-Null\n        $d.Bookmarks.Add(\"_GoBack\", $r)\n    }\n}\n"}
